# add disable property for IObject
#
# Adds a new "Disable" column (M) to the IObject sheet, mirroring the
# existing "MasterID" column (L) for formatting/values, and extends the
# sheet's data validations to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy column L (MasterID) formatting onto the new column M -------
$ws.Range("L1:L10").Copy()
$ws.Range("M1:M10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. New column header + values --------------------------------------
$ws.Range("M1").Value = "Disable"
$ws.Range("M2").Value = "int"
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("M5").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("M10").ClearContents()

# --- 3. Extend data validations to include column M ---------------------
$ws.Cells.Validation.Delete()

$ws.Range("F1").Validation.Add(0, 1, 1, "")
$ws.Range("I1:M1").Validation.Add(0, 1, 1, "")
$ws.Range("A7:A9").Validation.Add(0, 1, 1, "")

$ws.Range("F2").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("I2").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("B3:C3").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("F3:M6").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("F10:F1048576").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("I10:M1048576").Validation.Add(3, 1, "", """TRUE,FALSE""")
$ws.Range("B7:M9").Validation.Add(3, 1, "", """TRUE,FALSE""")

# --- 4. Move the active selection to the new header cell -----------------
$ws.Range("M1").Select()
